$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-7: Starting/Ending SoC (%) values swap
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 27

# Rows 8-43: relabel metrics with units and shift values (new row 43 added)
$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("B8").Value = 40.59608462201849
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("B9").Value = 34.11009726441294
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 72
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("B12").Value = 6239.258608
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("B13").Value = -1290.795189549715
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B14").Value = 104.6671430580555
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 7.027453633803263
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.466
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.077
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("B18").Value = 0.3890000000000002
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("B19").Value = 35
$ws.Range("A20").Value = "Maximum Temperature(C)"
$ws.Range("B20").Value = 47
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 12
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B22").Value = 71
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B23").Value = 67
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B24").Value = 67
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("B25").Value = 47
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("B26").Value = 0
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B27").Value = 0
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("B28").Value = 47
$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("B29").Value = 35
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B30").Value = 12
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 56
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.51206286
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001077796923559433
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 4.500972235504128
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 8.112588059035415
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 8.380351279844442
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 12.06528322335915
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 17.33766854738453
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 17.84769372987791
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 23.9361193458927
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 7.586624589589111
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0.0637531478116732
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
